$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was collected; insert it as a new row 63,
# pushing the existing row 63 (and everything after it) down by one.
$ws.Rows("63:63").Insert()

$ws.Range("A63").Value = 7
$ws.Range("B63").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C63").Value = "Ñuble"
$ws.Range("D63").Value = 44495
$ws.Range("E63").Value = 16
$ws.Range("F63").Value = 100112043
$ws.Range("G63").Value = "Pepino ensalada"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 100
$ws.Range("K63").Value = 13000
$ws.Range("L63").Value = 14000
$ws.Range("M63").Value = 13500
$ws.Range("N63").Value = "$/caja 60 unidades"
$ws.Range("O63").Value = "Región de Arica y Parinacota"
$ws.Range("P63").Value = 225
$ws.Range("Q63").Value = 60
$ws.Range("R63").Value = "Hortaliza"
